$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.395.71"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "3.209.96"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.596"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.59%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -2.40%  "
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "3.763.77"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D14").Value = "65.379.39"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").Value = "3.199.19"
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "413.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.203"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("E26").Value = "  -5.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("E33").Value = "  -1.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "157.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "2.734.73"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.717"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0635"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.15%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "298.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0263"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  -2.30%  "
$ws.Range("E48").Value = "  -8.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.909"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.55%  "
